$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-27: new ticker symbols per column ---
# Row 2
$ws.Range("B2").Value = "NSE:ASTRAMICRO"
$ws.Range("C2").Value = "NSE:APTECHT"
$ws.Range("D2").Value = "NSE:BHARTIARTL"
$ws.Range("E2").Value = "NSE:ESCORTS"
$ws.Range("F2").Value = "NSE:INDUSINDBK"

# Row 3
$ws.Range("B3").Value = "NSE:BSLGOLDETF"
$ws.Range("C3").Value = "NSE:BBL"
$ws.Range("D3").Value = "NSE:INDUSINDBK"
$ws.Range("E3").Value = "NSE:GAIL"
$ws.Range("F3").Value = ""

# Row 4
$ws.Range("B4").Value = "NSE:CANTABIL"
$ws.Range("C4").Value = "NSE:BERGEPAINT"
$ws.Range("D4").Value = "NSE:NAUKRI"
$ws.Range("E4").Value = "NSE:PFC"
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("B5").Value = "NSE:CENTURYPLY"
$ws.Range("C5").Value = "NSE:DELTAMAGNT"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "NSE:PIIND"
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("B6").Value = "NSE:DCM"
$ws.Range("C6").Value = "NSE:DHANI"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# Row 7
$ws.Range("B7").Value = "NSE:GREENLAM"
$ws.Range("C7").Value = "NSE:DNAMEDIA"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

# Row 8
$ws.Range("B8").Value = "NSE:KAYNES"
$ws.Range("C8").Value = "NSE:DUCON"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("B9").Value = "NSE:MITTAL"
$ws.Range("C9").Value = "NSE:ENGINERSIN"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("B10").Value = "NSE:MOHITIND"
$ws.Range("C10").Value = "NSE:FCSSOFT"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# Row 11
$ws.Range("B11").Value = "NSE:POWERINDIA"
$ws.Range("C11").Value = "NSE:FIEMIND"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# Row 12
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:GAIL"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

# Row 13
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:GENSOL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""

# Row 14
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:GOYALALUM"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""

# Row 15
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:GREENPOWER"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""

# Row 16
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:GULFOILLUB"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""

# Row 17
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:HDFCPVTBAN"
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""

# Row 18
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:HPIL"
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""

# Row 19
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "NSE:KABRAEXTRU"
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""

# Row 20
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "NSE:KIRLOSENG"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""

# Row 21
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = "NSE:LIBAS"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""

# Row 22
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "NSE:MACPOWER"
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""

# Row 23
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = "NSE:MAHABANK"
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = ""

# Row 24
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "NSE:MAPMYINDIA"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = ""

# Row 25
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = "NSE:MBLINFRA"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""

# Row 26
$ws.Range("B26").Value = ""
$ws.Range("C26").Value = "NSE:MIDHANI"
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""

# Row 27
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "NSE:MPSLTD"
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = ""

# --- Append new rows 28-35, extending the sheet dimension to A1:F35 ---
# Clone the bold/bordered/centered style used by the existing index column (A2:A27)
$ws.Range("A2").Copy()
$ws.Range("A28:A35").PasteSpecial(-4122)
for ($row = 28; $row -le 35; $row++) {
    $ws.Range("A$row").Value = $row - 2
}

# New rows only populate column C (support Zone); B, D, E, F stay blank
$ws.Range("C28").Value = "NSE:NESCO"
$ws.Range("C29").Value = "NSE:NOCIL"
$ws.Range("C30").Value = "NSE:OBCL"
$ws.Range("C31").Value = "NSE:PDMJEPAPER"
$ws.Range("C32").Value = "NSE:PGIL"
$ws.Range("C33").Value = "NSE:POKARNA"
$ws.Range("C34").Value = "NSE:RAJESHEXPO"
$ws.Range("C35").Value = "NSE:SAMBHAAV"
